# Big Update - HIV Surveillance manager tracks source, urls, details of data
# Adds a "Source" / "URL" footer block to the Stratified_Data sheet and tidies
# up leftover multi-pane selections on the Comments sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Stratified_Data: append two new rows under the existing table describing
# where the surveillance numbers came from.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Stratified_Data")

$ws1.Range("A19").Value = "Source"
$ws1.Range("B19:Z19").Value = "Florida Health Department"

$ws1.Range("A20").Value = "URL"
$ws1.Range("B20:Z20").Value = "testurl.org"

# Match the author's last on-screen selection (the newly added URL row).
$ws1.Range("C20:Z20").Select()

# ---------------------------------------------------------------------------
# Comments: the sheet was left with a stale split/multi-cell selection --
# clear it back down to a plain single view.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Comments")
$ws3.Range("A1").Select()

# Restore Stratified_Data as the active/visible tab (matches the original
# workbook state -- Comments should not stay selected).
$ws1.Activate()
$ws1.Range("C20:Z20").Select()
